$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-organization"
$wsMetadata.Range("B3").Value = "8.0.0"
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# Row 2 (Extension) Constraint(s) column is cleared
$wsElements.Range("AI2").Value = ""
# Row 5 (Extension.url) Fixed Value mirrors the new canonical URL
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-organization"
